# "sb sample name change"
#
# 1. The "No. of samples" label in the Grading block (cell A17) is renamed
#    to "No. of samples_scl" so it no longer collides with the identical
#    label used further down in the Nutritional Analysis block (A20/A22/
#    A24/A26, which stays "No. of samples ").
# 2. The four "No. of samples " rows in the Nutritional Analysis block
#    (rows 20, 22, 24 and 26) get their previously-blank B/D/F cells
#    filled in with explicit 0s, and the whole row's shading is unified
#    to match the rest of the sheet (light grey fill instead of white).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the sample-count label in the Grading block ---------------
$ws.Range("A17").Value = "No. of samples_scl"

# --- 2. Fix up the four "No. of samples " rows below -----------------
$lightGrey = 228 + 228 * 256 + 228 * 65536   # RGB(228,228,228) = #E4E4E4

foreach ($r in 20, 22, 24, 26) {
    $ws.Range("B$r").Value = 0
    $ws.Range("D$r").Value = 0
    $ws.Range("F$r").Value = 0

    foreach ($col in "B", "C", "D", "E", "F", "G") {
        $ws.Range("$col$r").Interior.Color = $lightGrey
    }
}

# --- cosmetic: restore the view state (scrolled down, K17 selected) ------
$ws.Range("K17").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
